# Fruta / hortaliza, semanal
#
# The weekly refresh adds one new "Frambuesa" price-report record for the
# Mercado Mayorista Lo Valledor de Santiago sheet. The new record is
# inserted as row 158 (pushing the existing rows 158-186 down to 159-187),
# which is why the sheet's used range grows from A1:T186 to A1:T187.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 158 - this shifts every row that was
# at 158..186 down by one (to 159..187) and extends the sheet dimension.
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(158, 1).Value  = 6
$ws.Cells.Item(158, 2).Value  = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(158, 3).Value  = 'Metropolitana'
$ws.Cells.Item(158, 4).Value  = 44642
$ws.Cells.Item(158, 5).Value  = 13
$ws.Cells.Item(158, 6).Value  = 'Fruta'
$ws.Cells.Item(158, 7).Value  = 100101
$ws.Cells.Item(158, 8).Value  = 'Berries'
$ws.Cells.Item(158, 9).Value  = 100101004
$ws.Cells.Item(158, 10).Value = 'Frambuesa'
$ws.Cells.Item(158, 11).Value = 'Sin especificar'
$ws.Cells.Item(158, 12).Value = 'Especial'
$ws.Cells.Item(158, 13).Value = 250
$ws.Cells.Item(158, 14).Value = 8000
$ws.Cells.Item(158, 15).Value = 8000
$ws.Cells.Item(158, 16).Value = 8000
$ws.Cells.Item(158, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(158, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(158, 19).Value = 4000
$ws.Cells.Item(158, 20).Value = 2
